$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing row-1 data before it gets overwritten by the new header.
$oldA1 = $ws.Range("A1").Value2
$oldB1 = $ws.Range("B1").Value2

# Copy the current A1 formatting (bold font, thin border, centered) so it can
# be re-applied to the new header cells after they are written.
$ws.Range("A1").Copy() | Out-Null

# Push the original data down to row 2 (A2 keeps it as a plain number, B2 as
# plain text - neither should carry any special style).
$ws.Range("A2").Value2 = $oldA1
$ws.Range("B2").Value2 = $oldB1

# Write the new header row.
$ws.Range("A1").Value2 = "ID"
$ws.Range("B1").Value2 = "name"

# Re-apply the header style (captured above) to both header cells.
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
